# Commit: "Fruta / hortaliza, semanal" - weekly price update.
# A new weekly price record (Ajo / Chino, Feria Lagunitas de Puerto Montt)
# is inserted at row 366, pushing the existing rows 366-442 down to 367-443.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 366 - shifts rows 366:442 down to 367:443.
$ws.Rows.Item(366).Insert()

# Populate the new row with the latest weekly data point.
$ws.Cells.Item(366, 1).Value = 4
$ws.Cells.Item(366, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(366, 3).Value = "Los Lagos"
$ws.Cells.Item(366, 4).Value = 45015
$ws.Cells.Item(366, 5).Value = 10
$ws.Cells.Item(366, 6).Value = 100112003
$ws.Cells.Item(366, 7).Value = "Ajo"
$ws.Cells.Item(366, 8).Value = "Chino"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 80
$ws.Cells.Item(366, 11).Value = 21000
$ws.Cells.Item(366, 12).Value = 21000
$ws.Cells.Item(366, 13).Value = 21000
$ws.Cells.Item(366, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(366, 15).Value = "China"
$ws.Cells.Item(366, 16).Value = 2100
$ws.Cells.Item(366, 17).Value = 10
$ws.Cells.Item(366, 18).Value = "Hortaliza"
